$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.255.94"
$ws.Range("E2").Value = "  +0.51%  "

$ws.Range("D3").Value = "1.859.10"
$ws.Range("E3").Value = "  +0.76%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9998"
$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7000"
$ws.Range("E5").Value = "  -0.11%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "237.79"
$ws.Range("E6").Value = "  +0.15%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.001"
$ws.Range("E7").Value = "  +0.12%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07996"
$ws.Range("E8").Value = "  +7.60%  "

$ws.Range("E9").Value = "  +0.66%  "

$ws.Range("E10").Value = "  -0.24%  "

$ws.Range("E11").Value = "  +0.97%  "

$ws.Range("D12").Value = "1.851.09"
$ws.Range("E12").Value = "  +0.24%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.7181"
$ws.Range("E13").Value = "  -0.90%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.178"
$ws.Range("E14").Value = "  -0.65%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "89.07"
$ws.Range("E15").Value = "  +0.05%  "

$ws.Range("D16").Value = "29.256.48"
$ws.Range("E16").Value = "  +0.42%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.761"
$ws.Range("E17").Value = "  -0.45%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.34"
$ws.Range("E18").Value = "  +2.74%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007791"
$ws.Range("E19").Value = "  +1.58%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "236.87"
$ws.Range("E20").Value = "  -1.99%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.001"
$ws.Range("E21").Value = "  +0.11%  "

$ws.Range("D22").Value = "2.108.13"
$ws.Range("E22").Value = "  +0.68%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.002"
$ws.Range("E23").Value = "  +0.20%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.454"
$ws.Range("E24").Value = "  -1.59%  "

$ws.Range("E25").Value = "  -0.22%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.993"
$ws.Range("E26").Value = "  +0.58%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1456"
$ws.Range("E27").Value = "  -1.06%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.07"
$ws.Range("E28").Value = "  +0.18%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.005"
$ws.Range("E29").Value = "  +3.64%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.437"
$ws.Range("E30").Value = "  +4.55%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.414"
$ws.Range("E31").Value = "  -1.04%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.484"
$ws.Range("E32").Value = "  -0.32%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.051"
$ws.Range("E33").Value = "  +1.11%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05220"
$ws.Range("E34").Value = "  +0.35%  "

$ws.Range("E35").Value = "  -1.10%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7061"
$ws.Range("E36").Value = "  -0.68%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9999"
$ws.Range("E37").Value = "  -0.30%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.664"
$ws.Range("E38").Value = "  +0.64%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01842"
$ws.Range("E39").Value = "  -1.52%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.719"

$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").Value = "1.143.30"
$ws.Range("E41").Value = "  +8.93%  "

$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9225"
$ws.Range("E42").Value = "  +1.46%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.928"
$ws.Range("E43").Value = "  +0.35%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4278"
$ws.Range("E44").Value = "  -0.32%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "70.71"
$ws.Range("E45").Value = "  +0.92%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.001"
$ws.Range("E46").Value = "  +0.14%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "103.51"
$ws.Range("E47").Value = "  +1.92%  "

$ws.Range("E48").Value = "  +2.38%  "

$ws.Range("D49").Value = "2.006.27"
$ws.Range("E49").Value = "  +0.29%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.187"
$ws.Range("E50").Value = "  -0.01%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.979"
$ws.Range("E51").Value = "  -1.84%  "
